$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.212283611297607
$ws.Range("B1").Value = 2.576965093612671
$ws.Range("C1").Value = 4.307102680206299
$ws.Range("D1").Value = 2.036807775497437
$ws.Range("E1").Value = 1.167977809906006
